$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.658.39"
$ws.Range("E2").Value = "  +0.16%  "

$ws.Range("D3").Value = "1.597.31"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'211.77"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("E6").Value = "  +1.21%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  +0.47%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").Value = "'19.52"
$ws.Range("E10").Value = "  +0.09%  "

$ws.Range("E11").Value = "  +0.49%  "

$ws.Range("D12").Value = "1.821.09"
$ws.Range("E12").Value = "  +0.65%  "

$ws.Range("D13").Value = "1.596.65"
$ws.Range("E13").Value = "  +0.36%  "

$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").Value = "'64.48"
$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("D17").Value = "26.640.30"
$ws.Range("E17").Value = "  +0.18%  "

$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  +0.73%  "

$ws.Range("D19").Value = "'208.80"
$ws.Range("E19").Value = "  -0.24%  "

$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("E21").Value = "  +3.95%  "

$ws.Range("D22").Value = "'4.28"
$ws.Range("E22").Value = "  +0.55%  "

$ws.Range("E23").Value = "  -2.48%  "

$ws.Range("D24").Value = "'8.90"
$ws.Range("E24").Value = "  +0.58%  "

$ws.Range("D25").Value = "'145.38"
$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("E27").Value = "  -0.57%  "

$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("D30").Value = "'0.0508"
$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("E31").Value = "  +0.79%  "

$ws.Range("E32").Value = "  +0.48%  "

$ws.Range("D33").Value = "'0.658"
$ws.Range("E33").Value = "  -3.06%  "

$ws.Range("E34").Value = "  +0.83%  "

$ws.Range("D35").Value = "1.276.39"
$ws.Range("E35").Value = "  -2.61%  "

$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("E37").Value = "  +0.68%  "

$ws.Range("E38").Value = "  -0.35%  "

$ws.Range("E39").Value = "  +2.32%  "

$ws.Range("E40").Value = "  +0.10%  "

$ws.Range("E41").Value = "  +1.72%  "

$ws.Range("D42").Value = "'64.66"
$ws.Range("E42").Value = "  +3.31%  "

$ws.Range("E43").Value = "  +1.66%  "

$ws.Range("D44").Value = "'0.787"
$ws.Range("E44").Value = "  -0.43%  "

$ws.Range("D45").Value = "1.734.23"
$ws.Range("E45").Value = "  +0.67%  "

$ws.Range("D46").Value = "'0.910"
$ws.Range("E46").Value = "  +8.57%  "

$ws.Range("D47").Value = "'90.11"
$ws.Range("E47").Value = "  +0.85%  "

$ws.Range("E48").Value = "  +0.25%  "

$ws.Range("E49").Value = "  +4.62%  "

$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.47"
$ws.Range("E51").Value = "  -0.68%  "
